# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.490.62"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "2.530.69"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.43"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.13"
$ws.Range("E6").Value = "  -2.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("D9").Value = "2.536.44"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0994"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.39"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "2.978.08"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.19"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").Value = "59.409.40"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").Value = "2.536.14"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.98"
$ws.Range("E19").Value = "  -2.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.23"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.25"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.95"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.422"
$ws.Range("E25").Value = "  -3.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.84"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.75"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "0.0₃0772"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.80"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.30"
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("E34").Value = "  -4.74%  "
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.58"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  -3.34%  "
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.03"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.66"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.814"
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.76"
$ws.Range("E42").Value = "  -5.74%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.28"
$ws.Range("E43").Value = "  -6.38%  "
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.603"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.88"
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "124.47"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.59"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0511"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0223"
$ws.Range("E51").Value = "  -1.72%  "
